$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full updated LR-pair table (ECs/FAPs/MuSCs/Resolving-Mac cross join for Ptn-Ptprb), rows 2-17
$rows = @(
    ,@("ECs", "Ptn", "Ptprb", "ECs", 3, 1, 1.00581, 3.01743, 0.003799625168827527, 0.003799625168827527, 3, 1, 237.3377026666667, 712.013108, 0.9887685707142667, 0.9887685707142667, 238.71663471916, 2148.44971247244, 0.003756949947431548, 0.003756949947431548)
    ,@("ECs", "Ptn", "Ptprb", "FAPs", 3, 1, 1.00581, 3.01743, 0.003799625168827527, 0.003799625168827527, 3, 1, 1.025352333333333, 3.076057, 0.004271702935173513, 0.004271702935173513, 1.03130963039, 9.28178667351, 0.0000162308699862397, 0.0000162308699862397)
    ,@("ECs", "Ptn", "Ptprb", "MuSCs", 3, 1, 1.00581, 3.01743, 0.003799625168827527, 0.003799625168827527, 3, 1, 1.587950666666667, 4.763852, 0.006615534293133127, 0.006615534293133127, 1.59717666004, 14.37458994036, 0.00002513655060543025, 0.00002513655060543025)
    ,@("ECs", "Ptn", "Ptprb", "Resolving-Mac", 3, 1, 1.00581, 3.01743, 0.003799625168827527, 0.003799625168827527, 3, 1, 0.08261766666666666, 0.247853, 0.0003441920574266213, 0.0003441920574266213, 0.08309767530999999, 0.74787907779, 0.00000130780080430872, 0.00000130780080430872)
    ,@("FAPs", "Ptn", "Ptprb", "ECs", 3, 1, 253.9688163333334, 761.9064490000001, 0.9594121222074437, 0.9594121222074438, 3, 1, 237.3377026666667, 712.013108, 0.9887685707142667, 0.9887685707142667, 60276.37541752595, 542487.3787577335, 0.9486365528009955, 0.9486365528009956)
    ,@("FAPs", "Ptn", "Ptprb", "FAPs", 3, 1, 253.9688163333334, 761.9064490000001, 0.9594121222074437, 0.9594121222074438, 3, 1, 1.025352333333333, 3.076057, 0.004271702935173513, 0.004271702935173513, 260.4075184212882, 2343.667665791593, 0.004098323578474586, 0.004098323578474587)
    ,@("FAPs", "Ptn", "Ptprb", "MuSCs", 3, 1, 253.9688163333334, 761.9064490000001, 0.9594121222074437, 0.9594121222074438, 3, 1, 1.587950666666667, 4.763852, 0.006615534293133127, 0.006615534293133127, 403.289951209061, 3629.609560881548, 0.006347023795710974, 0.006347023795710975)
    ,@("FAPs", "Ptn", "Ptprb", "Resolving-Mac", 3, 1, 253.9688163333334, 761.9064490000001, 0.9594121222074437, 0.9594121222074438, 3, 1, 0.08261766666666666, 0.247853, 0.0003441920574266213, 0.0003441920574266213, 20.98231101155522, 188.840799103997, 0.000330222032262621, 0.0003302220322626211)
    ,@("MuSCs", "Ptn", "Ptprb", "ECs", 3, 1, 9.336668333333334, 28.010005, 0.03527091597053946, 0.03527091597053946, 3, 1, 237.3377026666667, 712.013108, 0.9887685707142667, 0.9887685707142667, 2215.943412793949, 19943.49071514554, 0.03487477317197331, 0.03487477317197331)
    ,@("MuSCs", "Ptn", "Ptprb", "FAPs", 3, 1, 9.336668333333334, 28.010005, 0.03527091597053946, 0.03527091597053946, 3, 1, 1.025352333333333, 3.076057, 0.004271702935173513, 0.004271702935173513, 9.573374661142779, 86.160371950285, 0.0001506668752776117, 0.0001506668752776117)
    ,@("MuSCs", "Ptn", "Ptprb", "MuSCs", 3, 1, 9.336668333333334, 28.010005, 0.03527091597053946, 0.03527091597053946, 3, 1, 1.587950666666667, 4.763852, 0.006615534293133127, 0.006615534293133127, 14.82616870436222, 133.43551833926, 0.0002333359541533207, 0.0002333359541533207)
    ,@("MuSCs", "Ptn", "Ptprb", "Resolving-Mac", 3, 1, 9.336668333333334, 28.010005, 0.03527091597053946, 0.03527091597053946, 3, 1, 0.08261766666666666, 0.247853, 0.0003441920574266213, 0.0003441920574266213, 0.7713737521405555, 6.942363769265, 0.00001213996913522145, 0.00001213996913522145)
    ,@("Resolving-Mac", "Ptn", "Ptprb", "ECs", 2, 0.6666666666666666, 0.4016586666666667, 1.204976, 0.001517336653189343, 0.001517336653189343, 3, 1, 237.3377026666667, 712.013108, 0.9887685707142667, 0.9887685707142667, 95.32874520282311, 857.9587068254081, 0.001500294793866395, 0.001500294793866395)
    ,@("Resolving-Mac", "Ptn", "Ptprb", "FAPs", 2, 0.6666666666666666, 0.4016586666666667, 1.204976, 0.001517336653189343, 0.001517336653189343, 3, 1, 1.025352333333333, 3.076057, 0.004271702935173513, 0.004271702935173513, 0.4118416510702222, 3.706574859632, 0.00000648161143507527, 0.00000648161143507527)
    ,@("Resolving-Mac", "Ptn", "Ptprb", "MuSCs", 2, 0.6666666666666666, 0.4016586666666667, 1.204976, 0.001517336653189343, 0.001517336653189343, 3, 1, 1.587950666666667, 4.763852, 0.006615534293133127, 0.006615534293133127, 0.6378141475057778, 5.740327327552, 0.00001003799266340194, 0.00001003799266340194)
    ,@("Resolving-Mac", "Ptn", "Ptprb", "Resolving-Mac", 2, 0.6666666666666666, 0.4016586666666667, 1.204976, 0.001517336653189343, 0.001517336653189343, 3, 1, 0.08261766666666666, 0.247853, 0.0003441920574266213, 0.0003441920574266213, 0.03318410183644444, 0.298656916528, 0.0000005222552244700636, 0.0000005222552244700636)
)

$r = 2
foreach ($row in $rows) {
    $c = 1
    foreach ($val in $row) {
        $ws.Cells.Item($r, $c).Value = $val
        $c = $c + 1
    }
    $r = $r + 1
}
